$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '62.801.47'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +5.25%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.108.90'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +3.08%  '
$ws.Range("E4").Value = '  +0.07%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '586.04'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +3.73%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '143.97'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +3.04%  '
$ws.Range("E7").Value = '  +0.01%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '3.100.66'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +2.96%  '
$ws.Range("E9").Value = '  +1.65%  '
$ws.Range("E10").Value = '  +10.97%  '
$ws.Range("E11").Value = '  +7.80%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.0000244'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +5.14%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '35.42'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +3.97%  '
$ws.Range("E15").Value = '  +0.00%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.621.76'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +2.97%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '7.20'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.15%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.103.87'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +2.83%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '62.737.12'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +5.18%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '463.33'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +6.11%  '
$ws.Range("E21").Value = '  +2.75%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.729'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +5.55%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '13.37'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +0.81%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '82.24'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.86%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("E29").Value = '  +0.06%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '8.27'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +5.51%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '6.81'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +8.18%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '0.111'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +8.50%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '26.93'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +3.40%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0824'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +4.90%  '
$ws.Range("E35").Value = '  +11.29%  '
$ws.Range("E36").Value = '  +3.49%  '
$ws.Range("E37").Value = '  +1.46%  '
$ws.Range("E38").Value = '  +12.66%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '50.99'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +3.69%  '
$ws.Range("E40").Value = '  +1.63%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '428.97'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +6.28%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '2.908.62'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +4.41%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.0368'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +3.73%  '
$ws.Range("E44").Value = '  +9.34%  '
$ws.Range("E45").Value = '  +2.88%  '
$ws.Range("E46").Value = '  +6.95%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '34.93'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +3.27%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '123.32'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("E50").Value = '  +0.51%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '24.68'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +4.73%  '
